$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: metadata "dimension/measure" type row
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "iaest-measure:sexo"
$ws.Range("G2").Value = "iaest-measure:regimen"
$ws.Range("H2").Value = "iaest-measure:direccion-provincial-nombre"
$ws.Range("I2").Value = "sdmx-dimension:refArea"

# Row 3: metadata "medida/dim" row
$ws.Range("C3").Value = "dim"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "medida"
$ws.Range("H3").Value = "medida"

# Row 4: metadata "type/uri" row
$ws.Range("C4").Value = "URI-Municipio"
$ws.Range("F4").Value = "xsd:int"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("H4").Value = "xsd:int"
$ws.Range("I4").Value = "URI-comarca"

# Row 5 no longer exists - remove the mapping file reference row entirely
$ws.Range("F5:G5").EntireRow.Delete()
